# The deck ships with two theme parts: the one actually driving the slide
# master ("Integral" - green/yellow palette) and a second, unused "Office
# Theme" palette that only the notes master points at. The edit swaps which
# palette is "live": the slide master's 12 theme colors change from the
# Integral palette to the standard Office Theme palette.
#
# PowerPoint's ColorScheme object exposes exactly those 12 theme colors, in
# clrScheme document order (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink),
# so we drive the swap through $p.SlideMaster.ColorScheme.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cs = $m.ColorScheme

function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme palette (the content that theme1.xml/theme2.xml should end
# up carrying on the live design), in clrScheme order.
$cs.Colors(1).RGB  = RGBVal 0x00 0x00 0x00   # dk1      000000
$cs.Colors(2).RGB  = RGBVal 0xFF 0xFF 0xFF   # lt1      FFFFFF
$cs.Colors(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2      44546A
$cs.Colors(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2      E7E6E6
$cs.Colors(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1  5B9BD5
$cs.Colors(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2  ED7D31
$cs.Colors(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3  A5A5A5
$cs.Colors(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4  FFC000
$cs.Colors(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5  4472C4
$cs.Colors(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6  70AD47
$cs.Colors(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink    0563C1
$cs.Colors(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink 954F72
